$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Propagate A1's bordered/centered header style down to A2:A3 (the row
# index column keeps that look in the new layout) before we touch A1 itself.
$ws.Range("A1").Copy()
$ws.Range("A2:A3").PasteSpecial(-4122)

# --- Header row (row 1) ---
# The old A1 header ("NPOZZO") is gone; the remaining headers are renamed.
$ws.Range("A1").Clear()
$ws.Range("B1").Value = "codice"
$ws.Range("C1").Value = "data"
$ws.Range("D1").Value = "val"

# --- Row 2 ---
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 704
$ws.Range("C2").Value = 40224
$ws.Range("D2").Value = 95.5

# --- Row 3 ---
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 704
$ws.Range("C3").Value = 40267
$ws.Range("D3").Value = 94.40000000000001

# Date-formatted column C for the two data rows: apply lowercase format
# first, then the final uppercase format (matches the target workbook's
# registered-but-unused numFmtId 164 alongside the applied numFmtId 165).
$ws.Range("C2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("C2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C3").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- Page margins: 0.7/0.7/0.75/0.75/0.3/0.3 in -> 0.75/0.75/1/1/0.5/0.5 in ---
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36
